$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.154.94"
$ws.Range("E2").Value = "  +6.02%  "
$ws.Range("D3").Value = "1.719.12"
$ws.Range("E3").Value = "  +3.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.25"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9970"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3696"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.54"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3353"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.188"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07483"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9974"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.278"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.07"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.927"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.67%  "
$ws.Range("D16").Value = "1.720.42"
$ws.Range("E16").Value = "  +3.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001079"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06664"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.11"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9977"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.39"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.088"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.01"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.10%  "
$ws.Range("D24").Value = "26.076.84"
$ws.Range("E24").Value = "  +5.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.472"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.492"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.16"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.31"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.319"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.80%  "
$ws.Range("D30").Value = "1.907.50"
$ws.Range("E30").Value = "  +3.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "129.31"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.101"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08523"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.720"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.96"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.372"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.284"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06221"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02291"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2136"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.540"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.56"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6187"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9980"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.832"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5897"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.72"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.018"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07278"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.20"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.70%  "
